# Auto-generated: update market-price derived columns (H..N) per scheduled Universalis refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2078.5557
$ws.Range("I40").Value = 1571.6666
$ws.Range("K40").Value = 1571.6666
$ws.Range("M40").Value = -1396.6666
$ws.Range("H58").Value = 1508.4
$ws.Range("I58").Value = 681
$ws.Range("J58").Value = 2749.5
$ws.Range("K58").Value = 2043
$ws.Range("L58").Value = 8248.5
$ws.Range("M58").Value = -1893
$ws.Range("N58").Value = -8548.5
$ws.Range("H86").Value = 3748.7778
$ws.Range("I86").Value = 2941.923
$ws.Range("K86").Value = 2941.923
$ws.Range("M86").Value = -1818.923
$ws.Range("H89").Value = 3748.7778
$ws.Range("I89").Value = 2941.923
$ws.Range("K89").Value = 14709.615
$ws.Range("M89").Value = -9093.614999999998
$ws.Range("H92").Value = 3450.0588
$ws.Range("I92").Value = 3446.9375
$ws.Range("K92").Value = 3446.9375
$ws.Range("M92").Value = -2198.9375
$ws.Range("H121").Value = 2901.75
$ws.Range("J121").Value = 2901.75
$ws.Range("L121").Value = 8705.25
$ws.Range("N121").Value = -12199.25
$ws.Range("H127").Value = 2384.7778
$ws.Range("I127").Value = 1077.1666
$ws.Range("K127").Value = 3231.4998
$ws.Range("M127").Value = 1728.5002
$ws.Range("H129").Value = 214287.94
$ws.Range("I129").Value = 1022.7857
$ws.Range("J129").Value = 3200000
$ws.Range("K129").Value = 3068.3571
$ws.Range("L129").Value = 9600000
$ws.Range("M129").Value = 1931.6429
$ws.Range("N129").Value = -9610000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1288.093
$ws.Range("I2").Value = 1204.5
$ws.Range("K2").Value = 1204.5
$ws.Range("M2").Value = -1091.5
$ws.Range("H32").Value = 7734.725
$ws.Range("I32").Value = 2987.2188
$ws.Range("K32").Value = 2987.2188
$ws.Range("M32").Value = -2700.2188
$ws.Range("H74").Value = 5194.3335
$ws.Range("I74").Value = 4650.8
$ws.Range("J74").Value = 5873.75
$ws.Range("K74").Value = 4650.8
$ws.Range("L74").Value = 5873.75
$ws.Range("M74").Value = -3776.8
$ws.Range("N74").Value = -7621.75
$ws.Range("H77").Value = 5194.3335
$ws.Range("I77").Value = 4650.8
$ws.Range("J77").Value = 5873.75
$ws.Range("K77").Value = 23254
$ws.Range("L77").Value = 29368.75
$ws.Range("M77").Value = -18886
$ws.Range("N77").Value = -38104.75
$ws.Range("H98").Value = 72538.30499999999
$ws.Range("J98").Value = 72538.30499999999
$ws.Range("L98").Value = 72538.30499999999
$ws.Range("N98").Value = -78528.30499999999
$ws.Range("H116").Value = 1288.093
$ws.Range("I116").Value = 1204.5
$ws.Range("K116").Value = 1204.5
$ws.Range("M116").Value = 1089.5
$ws.Range("H131").Value = 89090.91
$ws.Range("J131").Value = 89090.91
$ws.Range("L131").Value = 89090.91
$ws.Range("N131").Value = -99170.91
$ws.Range("H132").Value = 3113.138
$ws.Range("I132").Value = 2760.8462
$ws.Range("K132").Value = 8282.5386
$ws.Range("M132").Value = -5752.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1288.093
$ws.Range("I3").Value = 1204.5
$ws.Range("K3").Value = 1204.5
$ws.Range("M3").Value = -1090.5
$ws.Range("H86").Value = 2523.3
$ws.Range("I86").Value = 2463
$ws.Range("J86").Value = 2549.1428
$ws.Range("K86").Value = 2463
$ws.Range("L86").Value = 2549.1428
$ws.Range("M86").Value = -1340
$ws.Range("N86").Value = -4795.1428
$ws.Range("H89").Value = 2523.3
$ws.Range("I89").Value = 2463
$ws.Range("J89").Value = 2549.1428
$ws.Range("K89").Value = 12315
$ws.Range("L89").Value = 12745.714
$ws.Range("M89").Value = -6699
$ws.Range("N89").Value = -23977.714
$ws.Range("H99").Value = 7024
$ws.Range("I99").Value = 7024
$ws.Range("K99").Value = 7024
$ws.Range("M99").Value = -5526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 53208.133
$ws.Range("I16").Value = 12760.25
$ws.Range("K16").Value = 12760.25
$ws.Range("M16").Value = -12473.25
$ws.Range("H58").Value = 4438
$ws.Range("I58").Value = 4635.8887
$ws.Range("K58").Value = 4635.8887
$ws.Range("M58").Value = -4432.8887
$ws.Range("H99").Value = 6465.727
$ws.Range("I99").Value = 5384.1333
$ws.Range("K99").Value = 5384.1333
$ws.Range("M99").Value = -3886.1333
$ws.Range("H113").Value = 53208.133
$ws.Range("I113").Value = 12760.25
$ws.Range("K113").Value = 12760.25
$ws.Range("M113").Value = -10590.25
$ws.Range("H126").Value = 6465.727
$ws.Range("I126").Value = 5384.1333
$ws.Range("K126").Value = 16152.3999
$ws.Range("M126").Value = -13682.3999
$ws.Range("H136").Value = 4438
$ws.Range("I136").Value = 4635.8887
$ws.Range("K136").Value = 13907.6661
$ws.Range("M136").Value = -11357.6661
$ws.Range("H141").Value = 184458.83
$ws.Range("J141").Value = 215700
$ws.Range("L141").Value = 215700
$ws.Range("N141").Value = -226060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2449.4285
$ws.Range("J11").Value = 3387
$ws.Range("L11").Value = 10161
$ws.Range("N11").Value = -10441
$ws.Range("H33").Value = 122.888885
$ws.Range("J33").Value = 188
$ws.Range("L33").Value = 1128
$ws.Range("N33").Value = -1694
$ws.Range("H39").Value = 6849.625
$ws.Range("I39").Value = 3966.5
$ws.Range("J39").Value = 15499
$ws.Range("K39").Value = 11899.5
$ws.Range("L39").Value = 46497
$ws.Range("M39").Value = -11605.5
$ws.Range("N39").Value = -47085
$ws.Range("H80").Value = 4527.5
$ws.Range("J80").Value = 4527.5
$ws.Range("L80").Value = 13582.5
$ws.Range("N80").Value = -15454.5
$ws.Range("H83").Value = 4527.5
$ws.Range("J83").Value = 4527.5
$ws.Range("L83").Value = 40747.5
$ws.Range("N83").Value = -50107.5
$ws.Range("H103").Value = 338
$ws.Range("J103").Value = 428
$ws.Range("L103").Value = 1284
$ws.Range("N103").Value = -3042
$ws.Range("H136").Value = 6211.0967
$ws.Range("I136").Value = 5573.174
$ws.Range("K136").Value = 16719.522
$ws.Range("M136").Value = -11619.522

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 222844
$ws.Range("J95").Value = 222844
$ws.Range("L95").Value = 222844
$ws.Range("N95").Value = -228336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1632
$ws.Range("I22").Value = 1070.8334
$ws.Range("K22").Value = 1070.8334
$ws.Range("M22").Value = -775.8334
$ws.Range("H27").Value = 1632
$ws.Range("I27").Value = 1070.8334
$ws.Range("K27").Value = 1070.8334
$ws.Range("M27").Value = -963.8334
$ws.Range("H40").Value = 6523.5386
$ws.Range("I40").Value = 6464.48
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 6464.48
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -6328.48
$ws.Range("N40").Value = -8272
$ws.Range("H55").Value = 350.4
$ws.Range("I55").Value = 276.875
$ws.Range("J55").Value = 644.5
$ws.Range("K55").Value = 276.875
$ws.Range("L55").Value = 644.5
$ws.Range("M55").Value = -103.875
$ws.Range("N55").Value = -990.5
$ws.Range("H82").Value = 1922.25
$ws.Range("I82").Value = 2753.0833
$ws.Range("J82").Value = 676
$ws.Range("K82").Value = 2753.0833
$ws.Range("L82").Value = 676
$ws.Range("M82").Value = -2392.0833
$ws.Range("N82").Value = -1398
$ws.Range("H85").Value = 1922.25
$ws.Range("I85").Value = 2753.0833
$ws.Range("J85").Value = 676
$ws.Range("K85").Value = 2753.0833
$ws.Range("L85").Value = 676
$ws.Range("M85").Value = -1505.0833
$ws.Range("N85").Value = -3172
$ws.Range("H93").Value = 2250
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 2250
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 2250
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -4746
$ws.Range("H100").Value = 4500
$ws.Range("J100").Value = 4500
$ws.Range("L100").Value = 4500
$ws.Range("N100").Value = -5582

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1822.2
$ws.Range("I96").Value = 1337.6666
$ws.Range("J96").Value = 2549
$ws.Range("K96").Value = 1337.6666
$ws.Range("L96").Value = 2549
$ws.Range("M96").Value = 35.33339999999998
$ws.Range("N96").Value = -5295
$ws.Range("H100").Value = 1022.5
$ws.Range("I100").Value = 1074.7727
$ws.Range("J100").Value = 735
$ws.Range("K100").Value = 2149.5454
$ws.Range("L100").Value = 1470
$ws.Range("M100").Value = -1608.5454
$ws.Range("N100").Value = -2552
$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040
$ws.Range("H132").Value = 5258.6924
$ws.Range("I132").Value = 5469.4546
$ws.Range("K132").Value = 16408.3638
$ws.Range("M132").Value = -13878.3638
